# Apply updated cryptocurrency price/volume data (GitHub Actions scheduled refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '25.948.69'
$ws.Range('E2').Value = '  -0.79%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.746.00'
$ws.Range('E3').Value = '  -0.28%  '
$ws.Range('E4').Value = '  +0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '248.78'
$ws.Range('E5').Value = '  +5.17%  '
$ws.Range('E6').Value = '  -0.02%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.5053'
$ws.Range('E7').Value = '  -4.87%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2751'
$ws.Range('E8').Value = '  -2.07%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06184'
$ws.Range('E9').Value = '  +0.13%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07263'
$ws.Range('E10').Value = '  +1.22%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.742.73'
$ws.Range('E11').Value = '  -0.42%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.6547'
$ws.Range('E12').Value = '  +1.09%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '15.20'
$ws.Range('E13').Value = '  -1.90%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '4.652'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '77.71'
$ws.Range('E15').Value = '  -1.04%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.9998'
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '25.975.33'
$ws.Range('E18').Value = '  -0.25%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.85'
$ws.Range('E19').Value = '  +0.69%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '0.000006849'
$ws.Range('E20').Value = '  +1.18%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '1.967.32'
$ws.Range('E21').Value = '  -0.16%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.449'
$ws.Range('E22').Value = '  +2.34%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '8.728'
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '5.397'
$ws.Range('E24').Value = '  +2.89%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '136.63'
$ws.Range('E25').Value = '  -1.84%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.507'
$ws.Range('E26').Value = '  -1.04%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '15.25'
$ws.Range('E27').Value = '  -0.23%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.785'
$ws.Range('E28').Value = '  -0.92%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '105.63'
$ws.Range('E29').Value = '  +0.59%  '
$ws.Range('E30').Value = '  +2.44%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.08211'
$ws.Range('E31').Value = '  -1.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '3.649'
$ws.Range('E32').Value = '  +0.41%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.04683'
$ws.Range('E33').Value = '  +0.97%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.655'
$ws.Range('E34').Value = '  +0.43%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.9961'
$ws.Range('E35').Value = '  -1.20%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.6193'
$ws.Range('E36').Value = '  -1.90%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.757'
$ws.Range('E37').Value = '  +1.74%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01612'
$ws.Range('E38').Value = '  -0.53%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.929'
$ws.Range('E39').Value = '  -2.04%  '
$ws.Range('E40').Value = '  -0.03%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '100.59'
$ws.Range('E41').Value = '  -1.43%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.3935'
$ws.Range('E42').Value = '  +0.06%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.7607'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '5.007'
$ws.Range('E44').Value = '  -1.59%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.1153'
$ws.Range('E45').Value = '  -0.11%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '6.325'
$ws.Range('E46').Value = '  -0.58%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '55.70'
$ws.Range('E47').Value = '  +1.89%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.05281'
$ws.Range('E48').Value = '  -1.14%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '30.69'
$ws.Range('E49').Value = '  -1.00%  '
$ws.Range('B50').Value = 'EnergySwap'
$ws.Range('C50').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.553'
$ws.Range('E50').Value = '  -0.62%  '
$ws.Range('B51').Value = 'Decentraland'
$ws.Range('C51').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.3439'
$ws.Range('E51').Value = '  -1.35%  '
